$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunControl")
$ws.Rows("38:41").Insert()
$ws.Range("AC39").Value = "simple"
